$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update of league bases, 17-05-2024 13:59
# Rows 64/65, 108/109/110, 113/114 are match rows whose underlying
# records were reordered (swapped / rotated) between adjacent row
# positions; row 217-224 odds were refreshed in-place for upcoming
# fixtures. Every changed cell is written explicitly below.

$ws.Range("B64").Value = 7082624
$ws.Range("E64").Value = "Colo Colo"
$ws.Range("F64").Value = "Deportes Copiapo"
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 1
$ws.Range("J64").Value = 1.333
$ws.Range("K64").Value = 5
$ws.Range("L64").Value = 8
$ws.Range("M64").Value = 1.45
$ws.Range("N64").Value = 4.75
$ws.Range("O64").Value = 7
$ws.Range("P64").Value = -1.25
$ws.Range("S64").Value = 3
$ws.Range("T64").Value = 1.875
$ws.Range("U64").Value = 1.925
$ws.Range("W64").Value = 3.75
$ws.Range("Y64").Value = -1
$ws.Range("Z64").Value = 0.875
$ws.Range("AA64").Value = -1
$ws.Range("AB64").Value = 0.925
$ws.Range("B65").Value = 7157967
$ws.Range("E65").Value = "Huachipato"
$ws.Range("F65").Value = "Palestino"
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 2
$ws.Range("J65").Value = 2.375
$ws.Range("K65").Value = 3.2
$ws.Range("L65").Value = 3
$ws.Range("M65").Value = 2.75
$ws.Range("N65").Value = 3.2
$ws.Range("O65").Value = 2.7
$ws.Range("P65").Value = 0
$ws.Range("S65").Value = 2.5
$ws.Range("T65").Value = 2
$ws.Range("U65").Value = 1.8
$ws.Range("W65").Value = 2.2
$ws.Range("Y65").Value = 0
$ws.Range("Z65").Value = 0
$ws.Range("AA65").Value = 1
$ws.Range("AB65").Value = -1
$ws.Range("B108").Value = 6077498
$ws.Range("E108").Value = "Universidad Catolica"
$ws.Range("F108").Value = "Deportes Copiapo"
$ws.Range("G108").Value = 2
$ws.Range("I108").Value = "D"
$ws.Range("J108").Value = 1.65
$ws.Range("K108").Value = 3.8
$ws.Range("L108").Value = 5.25
$ws.Range("M108").Value = 1.909
$ws.Range("N108").Value = 3.6
$ws.Range("O108").Value = 4.2
$ws.Range("P108").Value = -0.5
$ws.Range("Q108").Value = 1.85
$ws.Range("R108").Value = 2
$ws.Range("S108").Value = 2.75
$ws.Range("T108").Value = 2.025
$ws.Range("U108").Value = 1.825
$ws.Range("W108").Value = 2.6
$ws.Range("X108").Value = -1
$ws.Range("Z108").Value = 1
$ws.Range("AA108").Value = 1.025
$ws.Range("AB108").Value = -1
$ws.Range("B109").Value = 6078265
$ws.Range("E109").Value = "Audax Italiano"
$ws.Range("F109").Value = "Magallanes"
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 2
$ws.Range("I109").Value = "A"
$ws.Range("J109").Value = 1.666
$ws.Range("K109").Value = 3.75
$ws.Range("L109").Value = 5
$ws.Range("M109").Value = 2.25
$ws.Range("N109").Value = 3.3
$ws.Range("O109").Value = 3.3
$ws.Range("P109").Value = -0.25
$ws.Range("Q109").Value = 1.95
$ws.Range("R109").Value = 1.85
$ws.Range("S109").Value = 2.5
$ws.Range("T109").Value = 1.8
$ws.Range("U109").Value = 2
$ws.Range("V109").Value = -1
$ws.Range("X109").Value = 2.3
$ws.Range("Y109").Value = -1
$ws.Range("Z109").Value = 0.8500000000000001
$ws.Range("AA109").Value = -1
$ws.Range("AB109").Value = 1
$ws.Range("B110").Value = 6078266
$ws.Range("E110").Value = "Palestino"
$ws.Range("F110").Value = "Curico Unido"
$ws.Range("G110").Value = 4
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = "H"
$ws.Range("J110").Value = 1.533
$ws.Range("K110").Value = 4
$ws.Range("L110").Value = 6
$ws.Range("M110").Value = 1.363
$ws.Range("N110").Value = 4.75
$ws.Range("O110").Value = 7.5
$ws.Range("P110").Value = -1.5
$ws.Range("Q110").Value = 2.025
$ws.Range("R110").Value = 1.825
$ws.Range("S110").Value = 3
$ws.Range("T110").Value = 1.9
$ws.Range("U110").Value = 1.95
$ws.Range("V110").Value = 0.363
$ws.Range("W110").Value = -1
$ws.Range("Y110").Value = 1.025
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 0.8999999999999999
$ws.Range("B113").Value = 6078996
$ws.Range("E113").Value = "Colo Colo"
$ws.Range("F113").Value = "Union Espanola"
$ws.Range("H113").Value = 2
$ws.Range("J113").Value = 1.4
$ws.Range("K113").Value = 4.333
$ws.Range("L113").Value = 7
$ws.Range("M113").Value = 1.285
$ws.Range("N113").Value = 5.5
$ws.Range("O113").Value = 11
$ws.Range("P113").Value = -1.5
$ws.Range("Q113").Value = 1.9
$ws.Range("R113").Value = 1.95
$ws.Range("S113").Value = 3
$ws.Range("T113").Value = 2
$ws.Range("U113").Value = 1.85
$ws.Range("X113").Value = 10
$ws.Range("Z113").Value = 0.95
$ws.Range("AB113").Value = 0.8500000000000001
$ws.Range("B114").Value = 6077767
$ws.Range("E114").Value = "Nublense"
$ws.Range("F114").Value = "Huachipato"
$ws.Range("H114").Value = 1
$ws.Range("J114").Value = 2.75
$ws.Range("K114").Value = 3.4
$ws.Range("L114").Value = 2.45
$ws.Range("M114").Value = 2.875
$ws.Range("N114").Value = 3.3
$ws.Range("O114").Value = 2.5
$ws.Range("P114").Value = 0
$ws.Range("Q114").Value = 2.05
$ws.Range("R114").Value = 1.8
$ws.Range("S114").Value = 2.25
$ws.Range("T114").Value = 1.8
$ws.Range("U114").Value = 2.05
$ws.Range("X114").Value = 1.5
$ws.Range("Z114").Value = 0.8
$ws.Range("AB114").Value = 1.05
$ws.Range("M217").Value = 1.55
$ws.Range("N217").Value = 4.1
$ws.Range("O217").Value = 5.5
$ws.Range("Q217").Value = 1.95
$ws.Range("R217").Value = 1.9
$ws.Range("T217").Value = 2
$ws.Range("U217").Value = 1.85
$ws.Range("M218").Value = 3
$ws.Range("O218").Value = 2.3
$ws.Range("Q218").Value = 1.85
$ws.Range("R218").Value = 2
$ws.Range("Q219").Value = 1.975
$ws.Range("R219").Value = 1.875
$ws.Range("M220").Value = 2.05
$ws.Range("N220").Value = 3.3
$ws.Range("O220").Value = 3.6
$ws.Range("P220").Value = -0.25
$ws.Range("Q220").Value = 1.8
$ws.Range("R220").Value = 2.05
$ws.Range("T220").Value = 1.825
$ws.Range("U220").Value = 2.025
$ws.Range("T221").Value = 1.875
$ws.Range("U221").Value = 1.975
$ws.Range("M223").Value = 2.9
$ws.Range("O223").Value = 2.35
$ws.Range("Q223").Value = 1.8
$ws.Range("R223").Value = 2.05
$ws.Range("M224").Value = 2.45
$ws.Range("O224").Value = 2.875
$ws.Range("Q224").Value = 1.775
$ws.Range("R224").Value = 2.1
